# 08001-伊藤忠商事.xlsx — add 2022-Q4 data.
#
# Net effect requested by the diff:
#  - A brand-new "2022-Q4" sheet is inserted right after "总计", with fresh
#    fund numbers.
#  - Every existing quarter sheet ("2022-Q3" .. "2021-Q1") shifts down by one
#    slot and simply inherits the data that used to belong to the sheet
#    before it (i.e. each quarter's numbers "age" into the next slot).
#  - The old last sheet "2020-Q4" is untouched content-wise, but a literal
#    duplicate of it becomes the new last tab (so the old sheet9 formatting
#    — bold style "s=1", tabSelected, tight page margins — travels with the
#    2020-Q4 data to the end, while the slot it vacates (now "2021-Q1")
#    switches to the regular formatting used by the other quarter sheets).
#  - The "总计" summary sheet gets a new row 2 (2022-Q4 / 0.04) and every
#    later row shifts down by one, including a brand new row 10 for the
#    2020-Q4 entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the current last sheet ("2020-Q4") so the new copy
# — placed at the very end — keeps its exact original data AND exact
# original formatting (style indices, tabSelected, tight page margins).
# The original sheet is renamed out of the way ("2021-Q1-staging") so the
# copy can take over the "2020-Q4" name; it is turned into the new
# "2021-Q1" sheet in the last step of the cascade below.
# ---------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)      # "2020-Q4" (original)
$lastSheet.Copy($null, $lastSheet)            # new copy placed right after it

$lastSheet.Name = "2021-Q1-staging"
$newLast = $wb.Worksheets.Item($count + 1)
$newLast.Name = "2020-Q4"

# ---------------------------------------------------------------------
# Step 2: cascade every quarter sheet's data one slot down — each sheet
# takes on the numbers that used to live one sheet earlier in tab order
# (renamed left-to-right so every target name is vacated just in time).
# ---------------------------------------------------------------------

function Set-QuarterRow($ws, $d, $e, $f, $g, $h) {
    $ws.Cells.Item(2, 4).NumberFormat = "@"
    $ws.Cells.Item(2, 4).Value = $d
    $ws.Cells.Item(2, 5).NumberFormat = "@"
    $ws.Cells.Item(2, 5).Value = $e
    $ws.Cells.Item(2, 6).NumberFormat = "@"
    $ws.Cells.Item(2, 6).Value = $f
    $ws.Cells.Item(2, 7).NumberFormat = "@"
    $ws.Cells.Item(2, 7).Value = $g
    $ws.Cells.Item(2, 8).Value = $h
}

# "2022-Q3"(old, pos2) -> "2022-Q4", takes brand-new 2022-Q4 data
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "2022-Q4"
Set-QuarterRow $ws2 "1.18" "92.77" "3.35" "0.0395" 5

# "2022-Q2"(old, pos3) -> "2022-Q3", takes old "2022-Q3" data
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "2022-Q3"
Set-QuarterRow $ws3 "1.12" "90.06" "2.86" "0.0320" 5

# "2022-Q1"(old, pos4) -> "2022-Q2", takes old "2022-Q2" data
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "2022-Q2"
Set-QuarterRow $ws4 "1.20" "88.32" "3.00" "0.0360" 5

# "2021-Q4"(old, pos5) -> "2022-Q1", takes old "2022-Q1" data
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "2022-Q1"
Set-QuarterRow $ws5 "1.35" "88.71" "3.73" "0.0504" 4

# "2021-Q3"(old, pos6) -> "2021-Q4", takes old "2021-Q4" data, and the D1
# header label flips from "基金金额" to "基金规模" (matching "2021-Q4"'s).
$ws6 = $wb.Worksheets.Item(6)
$ws6.Name = "2021-Q4"
Set-QuarterRow $ws6 "1.22" "90.04" "2.78" "0.0339" 5
$ws6.Cells.Item(1, 4).Value = "基金规模"

# "2021-Q2"(old, pos7) -> "2021-Q3", takes old "2021-Q3" data
$ws7 = $wb.Worksheets.Item(7)
$ws7.Name = "2021-Q3"
Set-QuarterRow $ws7 "1.22" "88.77" "2.61" "0.0318" 9

# "2021-Q1"(old, pos8) -> "2021-Q2", takes old "2021-Q2" data
$ws8 = $wb.Worksheets.Item(8)
$ws8.Name = "2021-Q2"
Set-QuarterRow $ws8 "1.34" "90.09" "2.57" "0.0344" 9

# "2020-Q4"(old, pos9, staged above) -> "2021-Q1", takes old "2021-Q1" data,
# and loses the "last sheet" look: style index matching the regular quarter
# sheets (copied from sheet2's already-correct header/A2 formatting), no
# tabSelected (handled globally by the final Activate() below), and normal
# (non-tight) page margins.
$ws9 = $lastSheet
$ws9.Name = "2021-Q1"
Set-QuarterRow $ws9 "1.35" "87.46" "3.45" "0.0466" 5

$ws2.Range("B1:H1").Copy()
$ws9.Range("B1:H1").PasteSpecial(-4122)
$ws2.Range("A2").Copy()
$ws9.Range("A2").PasteSpecial(-4122)

$ws9.PageSetup.LeftMargin = 54
$ws9.PageSetup.RightMargin = 54
$ws9.PageSetup.TopMargin = 72
$ws9.PageSetup.BottomMargin = 72
$ws9.PageSetup.HeaderMargin = 36
$ws9.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# Step 3: update "总计" — insert the new 2022-Q4 row, shift the rest down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows("2:2").Insert(-4121)

# New row 2's A cell needs the same bold/border style used by the rest of
# column A; copy it from row 3 (which just shifted down, still correctly
# styled) before filling in values.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").Style = "Normal"

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.04

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q3"
$summary.Cells.Item(3, 3).Value = 1
$summary.Cells.Item(3, 4).Value = 0.03

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q2"
$summary.Cells.Item(4, 3).Value = 1
$summary.Cells.Item(4, 4).Value = 0.04

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2022-Q1"
$summary.Cells.Item(5, 3).Value = 1
$summary.Cells.Item(5, 4).Value = 0.05

$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = "2021-Q4"
$summary.Cells.Item(6, 3).Value = 1
$summary.Cells.Item(6, 4).Value = 0.03

$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(7, 2).Value = "2021-Q3"
$summary.Cells.Item(7, 3).Value = 1
$summary.Cells.Item(7, 4).Value = 0.03

$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(8, 2).Value = "2021-Q2"
$summary.Cells.Item(8, 3).Value = 1
$summary.Cells.Item(8, 4).Value = 0.03

$summary.Cells.Item(9, 1).Value = 7
$summary.Cells.Item(9, 2).Value = "2021-Q1"
$summary.Cells.Item(9, 3).Value = 1
$summary.Cells.Item(9, 4).Value = 0.05

$summary.Cells.Item(10, 1).Value = 8
$summary.Cells.Item(10, 2).Value = "2020-Q4"
$summary.Cells.Item(10, 3).Value = 1
$summary.Cells.Item(10, 4).Value = 0.04

# Row 10's A cell needs the same bold/border style as the rest of column A.
$summary.Range("A9").Copy()
$summary.Range("A10").PasteSpecial(-4122)

# Make the final "2020-Q4" sheet the selected tab, matching the source tab
# order (it is the last, rightmost sheet again after the reshuffle). This
# also clears tabSelected from whichever sheet happened to carry it before.
$newLast.Activate()
